# Updates crypto price/volume figures to refreshed values from the data source.
# Cells store text-formatted numbers/percentages (inline strings), so we force
# a Text number format while writing the value to stop Excel auto-coercing the
# numeric-looking text into a real number, then restore the original cell style
# so formatting is unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Worksheet, $Address, $NewValue) {
    $cell = $Worksheet.Range($Address)
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = $originalStyle
}

Set-TextValue $ws "D2" "261.66"
Set-TextValue $ws "E2" "1.75%"
Set-TextValue $ws "D3" "27.23"
Set-TextValue $ws "D4" "4.722"
Set-TextValue $ws "E4" "3.07%"
Set-TextValue $ws "E5" "2.83%"
Set-TextValue $ws "D6" "6.638"
Set-TextValue $ws "E6" "0.17%"
Set-TextValue $ws "E7" "1.28%"
Set-TextValue $ws "D8" "0.9213"
Set-TextValue $ws "E8" "-2.46%"
Set-TextValue $ws "D9" "0.1408"
Set-TextValue $ws "E9" "1.37%"
Set-TextValue $ws "D10" "0.05102"
Set-TextValue $ws "E10" "3.32%"
Set-TextValue $ws "D11" "0.07112"
Set-TextValue $ws "E11" "0.53%"
Set-TextValue $ws "D12" "0.03043"
Set-TextValue $ws "E12" "-0.93%"
Set-TextValue $ws "D13" "0.09107"
Set-TextValue $ws "E13" "-0.29%"
Set-TextValue $ws "D14" "0.001531"
Set-TextValue $ws "E14" "-0.06%"
Set-TextValue $ws "D15" "0.0006115"
Set-TextValue $ws "E15" "0.88%"
Set-TextValue $ws "D16" "0.006152"
Set-TextValue $ws "E16" "2.02%"
Set-TextValue $ws "D17" "3.451"
Set-TextValue $ws "E17" "-1.29%"
Set-TextValue $ws "D18" "3.171"
Set-TextValue $ws "E18" "-0.37%"
Set-TextValue $ws "E20" "2.42%"
Set-TextValue $ws "E21" "2.19%"
Set-TextValue $ws "D22" "4.100"
Set-TextValue $ws "E22" "3.63%"
Set-TextValue $ws "D23" "0.04261"
Set-TextValue $ws "E23" "-0.43%"
Set-TextValue $ws "E24" "-0.39%"
Set-TextValue $ws "E25" "-8.78%"
Set-TextValue $ws "D26" "0.0001201"
Set-TextValue $ws "E26" "0.03%"
Set-TextValue $ws "E27" "3.11%"
Set-TextValue $ws "D40" "0.03881"
Set-TextValue $ws "E40" "1.58%"
Set-TextValue $ws "D41" "0.1113"
Set-TextValue $ws "E41" "1.05%"
Set-TextValue $ws "D42" "0.004130"
Set-TextValue $ws "E42" "6.05%"
Set-TextValue $ws "D43" "0.01485"
Set-TextValue $ws "E43" "4.77%"
Set-TextValue $ws "E44" "-9.91%"
Set-TextValue $ws "D45" "0.00005326"
Set-TextValue $ws "E45" "-0.99%"
Set-TextValue $ws "E47" "-18.44%"
